$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Metadata" (first sheet): bump Version/Date, fill Publisher, replace
# the duplicated "Contact" rows with "Jurisdiction" + the "Description" row
# that used to follow them, then delete the now-redundant old "Description"
# row so everything below shifts up by one (sheet shrinks from 21 to 20 rows).
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$meta.Range("B9").Value = "Alvearie Team"

$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

$meta.Range("A11").Value = "Description"
$meta.Range("B11").Value = "Count of medications for the episode of care"

# The old row 12 ("Description" / "Count of medications...") is now a
# duplicate of row 11 above, so remove it; every later row shifts up.
$meta.Rows.Item(12).Delete()

# ---------------------------------------------------------------------------
# Sheet "Elements" (second sheet): the root Extension row's Short/Definition
# columns (K2/L2) now surface the StructureDefinition's own Title/Description
# instead of the generic "Extension" / "An Extension".
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("K2").Value = "Episode Rx Count"
$elements.Range("L2").Value = "Count of medications for the episode of care"
